$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.904.46"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.584.45"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "521.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.40%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Value = "2.594.09"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.329"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("E13").Value = "  +3.18%  "
$ws.Range("D14").Value = "3.041.69"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "58.859.42"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.39"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "2.596.36"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  -1.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.403"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.43%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "148.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.824"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.809"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.24%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  +0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.01%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0951"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.588"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").Value = "1.969.71"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -0.45%  "
